$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 304; this shifts the existing rows 304-369
# down to 305-370 (matching the diff, which shows the whole block from
# row 304 onward being pushed down by one row after a new record is
# inserted).
$ws.Rows.Item(304).Insert()

# Populate the newly inserted row 304 with the new record's data.
$ws.Cells.Item(304, 1).Value2 = 6
$ws.Cells.Item(304, 2).Value2 = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(304, 3).Value2 = "Metropolitana"
$ws.Cells.Item(304, 4).Value2 = 45258
$ws.Cells.Item(304, 5).Value2 = 13
$ws.Cells.Item(304, 6).Value2 = 100112022
$ws.Cells.Item(304, 7).Value2 = "Arveja Verde"
$ws.Cells.Item(304, 8).Value2 = "Sin especificar"
$ws.Cells.Item(304, 9).Value2 = "Primera"
$ws.Cells.Item(304, 10).Value2 = 380
$ws.Cells.Item(304, 11).Value2 = 16000
$ws.Cells.Item(304, 12).Value2 = 18000
$ws.Cells.Item(304, 13).Value2 = 17053
$ws.Cells.Item(304, 14).Value2 = '$/saco 25 kilos'
$ws.Cells.Item(304, 15).Value2 = "Región del Maule"
$ws.Cells.Item(304, 16).Value2 = 682
$ws.Cells.Item(304, 17).Value2 = 25
$ws.Cells.Item(304, 18).Value2 = "Hortaliza"
